$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista de Itens de Trabalho")
$ws.Activate()

# Rows 46-55: mark "Estado Atual" (column D) as "Incompleto" instead of
# "Concluído", and reset "Horas Trabalhadas" (column H) to 0.
for ($r = 46; $r -le 55; $r++) {
    $ws.Cells.Item($r, 4).Value = "Incompleto"
    $ws.Cells.Item($r, 8).Value = 0
}

# Update the sheet view scroll position / selection to match the saved
# workbook state (top row 31 visible, H48 selected).
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H48").Select()
